# Auto-generated Excel COM-interop script
# Applies schedule update for Linea 141 across sheets LP1912, LP1912-215, 6203-6173
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2,1).Value = 'Última actualización: 07:58:19'
$ws.Cells.Item(3,1).Value = 'Total filas: 90'
$ws.Cells.Item(62,1).Value = '06:43:12'
$ws.Cells.Item(62,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(62,4).Value = 80
$ws.Cells.Item(63,1).Value = '07:19:37'
$ws.Cells.Item(63,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(63,4).Value = 44
$ws.Cells.Item(66,1).Value = '07:58:19'
$ws.Cells.Item(66,2).Value = '08:11'
$ws.Cells.Item(66,3).Value = '16_SANTA ANA'
$ws.Cells.Item(66,4).Value = 13
$ws.Cells.Item(67,1).Value = '06:14:19'
$ws.Cells.Item(67,2).Value = '08:12'
$ws.Cells.Item(67,3).Value = '15_ABASTO'
$ws.Cells.Item(67,4).Value = 118
$ws.Cells.Item(68,1).Value = '07:45:49'
$ws.Cells.Item(68,2).Value = '08:13'
$ws.Cells.Item(68,3).Value = '10_OLMOS'
$ws.Cells.Item(68,4).Value = 28
$ws.Cells.Item(69,1).Value = '06:43:12'
$ws.Cells.Item(69,2).Value = '08:21'
$ws.Cells.Item(69,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(69,4).Value = 98
$ws.Cells.Item(70,1).Value = '07:19:37'
$ws.Cells.Item(70,2).Value = '08:22'
$ws.Cells.Item(70,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(70,4).Value = 63
$ws.Cells.Item(71,3).Value = '215B_EL PATO'
$ws.Cells.Item(72,2).Value = '08:23'
$ws.Cells.Item(72,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(72,4).Value = 100
$ws.Cells.Item(73,1).Value = '06:43:12'
$ws.Cells.Item(73,2).Value = '08:27'
$ws.Cells.Item(73,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(73,4).Value = 104
$ws.Cells.Item(74,1).Value = '07:45:49'
$ws.Cells.Item(74,2).Value = '08:31'
$ws.Cells.Item(74,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(74,4).Value = 46
$ws.Cells.Item(75,1).Value = '07:58:19'
$ws.Cells.Item(75,2).Value = '08:33'
$ws.Cells.Item(75,3).Value = '10_OLMOS'
$ws.Cells.Item(75,4).Value = 35
$ws.Cells.Item(76,1).Value = '07:58:19'
$ws.Cells.Item(76,2).Value = '08:35'
$ws.Cells.Item(76,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(76,4).Value = 37
$ws.Cells.Item(77,1).Value = '06:57:11'
$ws.Cells.Item(77,2).Value = '08:42'
$ws.Cells.Item(77,3).Value = '81_EL PELIGRO'
$ws.Cells.Item(77,4).Value = 105
$ws.Cells.Item(78,1).Value = '07:58:19'
$ws.Cells.Item(78,2).Value = '08:42'
$ws.Cells.Item(78,3).Value = '16_SANTA ANA'
$ws.Cells.Item(78,4).Value = 44
$ws.Cells.Item(79,2).Value = '08:43'
$ws.Cells.Item(79,3).Value = '14_ABASTO'
$ws.Cells.Item(79,4).Value = 84
$ws.Cells.Item(80,1).Value = '06:57:11'
$ws.Cells.Item(80,2).Value = '08:54'
$ws.Cells.Item(80,3).Value = '17_ROMERO'
$ws.Cells.Item(81,1).Value = '07:19:37'
$ws.Cells.Item(81,2).Value = '09:01'
$ws.Cells.Item(81,3).Value = '215A_EL PATO'
$ws.Cells.Item(81,4).Value = 102
$ws.Cells.Item(82,2).Value = '09:03'
$ws.Cells.Item(82,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(82,4).Value = 78
$ws.Cells.Item(83,1).Value = '07:19:37'
$ws.Cells.Item(83,2).Value = '09:10'
$ws.Cells.Item(83,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(83,4).Value = 111
$ws.Cells.Item(84,1).Value = '07:19:37'
$ws.Cells.Item(84,2).Value = '09:16'
$ws.Cells.Item(84,3).Value = '27_EL RETIRO'
$ws.Cells.Item(84,4).Value = 117
$ws.Cells.Item(85,1).Value = '07:58:19'
$ws.Cells.Item(85,2).Value = '09:17'
$ws.Cells.Item(85,3).Value = '27_EL RETIRO'
$ws.Cells.Item(85,4).Value = 79
$ws.Cells.Item(86,2).Value = '09:21'
$ws.Cells.Item(86,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(86,4).Value = 96
$ws.Cells.Item(87,1).Value = '07:45:49'
$ws.Cells.Item(87,2).Value = '09:22'
$ws.Cells.Item(87,3).Value = '17_ROMERO'
$ws.Cells.Item(87,4).Value = 97
$ws.Cells.Item(87,5).Value = 'LP1912'
$ws.Cells.Item(88,1).Value = '07:45:49'
$ws.Cells.Item(88,2).Value = '09:23'
$ws.Cells.Item(88,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(88,4).Value = 98
$ws.Cells.Item(88,5).Value = 'LP1912'
$ws.Cells.Item(89,1).Value = '07:58:19'
$ws.Cells.Item(89,2).Value = '09:23'
$ws.Cells.Item(89,3).Value = '17_ROMERO'
$ws.Cells.Item(89,4).Value = 85
$ws.Cells.Item(89,5).Value = 'LP1912'
$ws.Cells.Item(90,1).Value = '07:45:49'
$ws.Cells.Item(90,2).Value = '09:32'
$ws.Cells.Item(90,3).Value = '15_ABASTO'
$ws.Cells.Item(90,4).Value = 107
$ws.Cells.Item(90,5).Value = 'LP1912'
$ws.Cells.Item(91,1).Value = '07:45:49'
$ws.Cells.Item(91,2).Value = '09:33'
$ws.Cells.Item(91,3).Value = '10_OLMOS'
$ws.Cells.Item(91,4).Value = 108
$ws.Cells.Item(91,5).Value = 'LP1912'
$ws.Cells.Item(92,1).Value = '07:45:49'
$ws.Cells.Item(92,2).Value = '09:41'
$ws.Cells.Item(92,3).Value = '215C_EL PATO'
$ws.Cells.Item(92,4).Value = 116
$ws.Cells.Item(92,5).Value = 'LP1912'
$ws.Cells.Item(93,1).Value = '07:58:19'
$ws.Cells.Item(93,2).Value = '09:42'
$ws.Cells.Item(93,3).Value = '215C_EL PATO'
$ws.Cells.Item(93,4).Value = 104
$ws.Cells.Item(93,5).Value = 'LP1912'
$ws.Cells.Item(94,1).Value = '07:58:19'
$ws.Cells.Item(94,2).Value = '09:43'
$ws.Cells.Item(94,3).Value = '14_ABASTO'
$ws.Cells.Item(94,4).Value = 105
$ws.Cells.Item(94,5).Value = 'LP1912'
$ws.Cells.Item(95,1).Value = '07:58:19'
$ws.Cells.Item(95,2).Value = '09:52'
$ws.Cells.Item(95,3).Value = '15_ABASTO'
$ws.Cells.Item(95,4).Value = 114
$ws.Cells.Item(95,5).Value = 'LP1912'

$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2,1).Value = 'Última actualización: 07:58:19'
$ws.Cells.Item(3,1).Value = 'Total filas: 14'
$ws.Cells.Item(19,1).Value = '07:58:19'
$ws.Cells.Item(19,2).Value = '09:42'
$ws.Cells.Item(19,3).Value = '215C_EL PATO'
$ws.Cells.Item(19,4).Value = 104
$ws.Cells.Item(19,5).Value = 'LP1912'

$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2,1).Value = 'Última actualización: 07:58:19'
$ws.Cells.Item(3,1).Value = 'Total filas: 19'
$ws.Cells.Item(19,1).Value = '07:58:19'
$ws.Cells.Item(19,2).Value = '08:11'
$ws.Cells.Item(19,3).Value = '215C_LA PLATA'
$ws.Cells.Item(19,4).Value = 13
$ws.Cells.Item(19,5).Value = 'L6203'
$ws.Cells.Item(20,1).Value = '07:19:37'
$ws.Cells.Item(20,2).Value = '08:35'
$ws.Cells.Item(20,4).Value = 76
$ws.Cells.Item(21,1).Value = '06:57:11'
$ws.Cells.Item(21,2).Value = '08:38'
$ws.Cells.Item(21,3).Value = '215A_LA PLATA'
$ws.Cells.Item(21,4).Value = 101
$ws.Cells.Item(21,5).Value = 'L6173'
$ws.Cells.Item(22,1).Value = '07:58:19'
$ws.Cells.Item(22,2).Value = '08:40'
$ws.Cells.Item(22,3).Value = '215A_LA PLATA'
$ws.Cells.Item(22,4).Value = 42
$ws.Cells.Item(22,5).Value = 'L6173'
$ws.Cells.Item(23,1).Value = '07:19:37'
$ws.Cells.Item(23,2).Value = '09:08'
$ws.Cells.Item(23,3).Value = '215D_LA PLATA'
$ws.Cells.Item(23,4).Value = 109
$ws.Cells.Item(23,5).Value = 'L6203'
$ws.Cells.Item(24,1).Value = '07:58:19'
$ws.Cells.Item(24,2).Value = '09:09'
$ws.Cells.Item(24,3).Value = '215D_LA PLATA'
$ws.Cells.Item(24,4).Value = 71
$ws.Cells.Item(24,5).Value = 'L6203'
